$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Set the text labels first, in the exact order the shared-string table
# needs them appended (49=UC01, 50=AD01 Toke, 51=Review UC02, 52=review AD02,
# 53=DOM01, 54=Review DOM02, 55=Pause).
$ws.Range("A11").Value = "UC01 med EMIL"
$ws.Range("A13").Value = "AD01 med Toke"
$ws.Range("A12").Value = "Review af UC02 med Emil"
$ws.Range("A14").Value = "review af AD02"
$ws.Range("A15").Value = "DOM01 med MIK"
$ws.Range("A16").Value = "Review af DOM02"
$ws.Range("A17").Value = "Pause "

# Row 11
$ws.Range("C11").Value = 43964
$ws.Range("D11").Value = 0.375
$ws.Range("E11").Value = 0.39583333333333331

# Row 12
$ws.Range("C12").Value = 43964
$ws.Range("D12").Value = 0.39583333333333331
$ws.Range("E12").Value = 0.41666666666666669

# Row 13
$ws.Range("C13").Value = 43964
$ws.Range("D13").Value = 0.41666666666666669
$ws.Range("E13").Value = 0.4375

# Row 14
$ws.Range("C14").Value = 43964
$ws.Range("D14").Value = 0.4375
$ws.Range("E14").Value = 0.45833333333333331

# Row 15
$ws.Range("C15").Value = 43964
$ws.Range("D15").Value = 0.45833333333333331
$ws.Range("E15").Value = 0.47916666666666669

# Row 16
$ws.Range("C16").Value = 43964
$ws.Range("D16").Value = 0.47916666666666669
$ws.Range("E16").Value = 0.5

# Row 17
$ws.Range("C17").Value = 43964
$ws.Range("D17").Value = 0.50694444444444442
$ws.Range("E17").Value = 0.54166666666666663

# Update the active selection to B14
$ws.Range("B14").Select()
